$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4551307559013367
$ws.Range("B1").Value = 3.342988967895508
$ws.Range("C1").Value = 4.272516250610352
$ws.Range("D1").Value = 1.572198271751404
$ws.Range("E1").Value = 1.195690512657166
